$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 20 data: Question No -> GFG, GFG/LC -> GFG, Question -> Bottom View of Binary Tree
$ws.Range("A20").Value = "GFG"
$ws.Range("B20").Value = "GFG"
$ws.Range("C20").Value = "Bottom View of Binary Tree"

# Match formatting of the existing rows (style index 1: left/top aligned, wrap text)
$ws.Range("A20:C20").HorizontalAlignment = -4131
$ws.Range("A20:C20").VerticalAlignment = -4160
$ws.Range("A20:C20").WrapText = $true

# Update selection to reflect the newly added row, as in the edited workbook
$ws.Range("A20:C20").Select()
